# Fruta / hortaliza, semanal
# Insert one new weekly record as row 4 (pushing the existing data rows
# 4-111 down to 5-112), matching the updated dimension A1:T112.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 4..111 down to 5..112, leaving row 4 free for the
# new record.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the latest weekly price observation.
$ws.Cells.Item(4, 1).Value = 9
$ws.Cells.Item(4, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(4, 3).Value = "Metropolitana"
$ws.Cells.Item(4, 4).Value = 44922
$ws.Cells.Item(4, 5).Value = 13
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100101
$ws.Cells.Item(4, 8).Value = "Berries"
$ws.Cells.Item(4, 9).Value = 100101004
$ws.Cells.Item(4, 10).Value = "Frambuesa"
$ws.Cells.Item(4, 11).Value = "Sin especificar"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 380
$ws.Cells.Item(4, 14).Value = 8000
$ws.Cells.Item(4, 15).Value = 8000
$ws.Cells.Item(4, 16).Value = 8000
$ws.Cells.Item(4, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(4, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(4, 19).Value = 4000
$ws.Cells.Item(4, 20).Value = 2
